# Update gh-pages to output generated at 456a3b4
# Updates attendee/follower counts (column F) across the four sheets:
#   展览 (Exhibitions), 演出 (Performances), 本地生活 (Local Life),
#   全部类型 (All Types - aggregate of the three sheets above)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5848
$ws.Range("F6").Value = 61
$ws.Range("F9").Value = 1563
$ws.Range("F11").Value = 27
$ws.Range("F12").Value = 665
$ws.Range("F15").Value = 1524
$ws.Range("F17").Value = 113
$ws.Range("F18").Value = 607
$ws.Range("F19").Value = 4345
$ws.Range("F20").Value = 25
$ws.Range("F22").Value = 3328
$ws.Range("F23").Value = 807
$ws.Range("F25").Value = 35
$ws.Range("F26").Value = 2284
$ws.Range("F27").Value = 44
$ws.Range("F30").Value = 449
$ws.Range("F31").Value = 1215
$ws.Range("F34").Value = 1181
$ws.Range("F35").Value = 1174

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 18
$ws.Range("F15").Value = 25
$ws.Range("F16").Value = 4
$ws.Range("F18").Value = 123
$ws.Range("F19").Value = 293
$ws.Range("F20").Value = 225
$ws.Range("F21").Value = 489

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 248
$ws.Range("F3").Value = 621
$ws.Range("F4").Value = 162
$ws.Range("F5").Value = 250

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 248
$ws.Range("F6").Value = 621
$ws.Range("F7").Value = 162
$ws.Range("F8").Value = 5848
$ws.Range("F11").Value = 61
$ws.Range("F15").Value = 18
$ws.Range("F20").Value = 1563
$ws.Range("F23").Value = 27
$ws.Range("F26").Value = 1524
$ws.Range("F28").Value = 113
$ws.Range("F29").Value = 607
$ws.Range("F30").Value = 4345
$ws.Range("F32").Value = 3328
$ws.Range("F33").Value = 807
$ws.Range("F34").Value = 2284
$ws.Range("F35").Value = 44
$ws.Range("F38").Value = 449
$ws.Range("F39").Value = 1215
$ws.Range("F41").Value = 123
$ws.Range("F42").Value = 293
$ws.Range("F43").Value = 225
$ws.Range("F44").Value = 489
$ws.Range("F49").Value = 1174
